$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 311, shifting existing rows 311:439 down to 312:440.
$ws.Rows.Item(311).Insert()

# Populate the newly inserted row 311 with the new weekly record.
# Columns A,B,C,E,F,G,H,I,N,Q,R carry over the same values as the
# (now shifted-down) neighbouring rows; D,J,K,L,M,O,P hold the new data.
$ws.Range("A311").Value = 3
$ws.Range("B311").Value = "Femacal de La Calera"
$ws.Range("C311").Value = "Coquimbo"
$ws.Range("D311").Value = 44784
$ws.Range("E311").Value = 5
$ws.Range("F311").Value = 100112017
$ws.Range("G311").Value = "Apio"
$ws.Range("H311").Value = "Americana (o)"
$ws.Range("I311").Value = "Primera"
$ws.Range("J311").Value = 230
$ws.Range("K311").Value = 9000
$ws.Range("L311").Value = 9500
$ws.Range("M311").Value = 9239
$ws.Range("N311").Value = "$/docena de matas"
$ws.Range("O311").Value = "Provincia de Limarí"
$ws.Range("P311").Value = 1540
$ws.Range("Q311").Value = 6
$ws.Range("R311").Value = "Hortaliza"
